$d = $word.ActiveDocument

# Locate the target paragraph: "Using cplex on the DTUHPC cluster:"
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*Using cplex on the*") {
        $target = $para
        break
    }
}

$pStart = $target.Range.Start

# 1) Remove the existing "_GoBack" bookmark (it currently sits in an empty
#    paragraph near the end of the document).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) Re-create the "_GoBack" bookmark, collapsed, at the very start of the
#    target paragraph -- this also acts as a run-boundary so later edits
#    don't get merged into the following runs.
$bmRange = $d.Range($pStart, $pStart)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 3) Insert the new " Tip " text right before the bookmark (i.e. at the very
#    start of the paragraph).
$insRange = $d.Range($pStart, $pStart)
$insRange.InsertBefore(" Tip ")

# 4) Make the newly inserted text bold.
$newRange = $d.Range($pStart, $pStart + 5)
$newRange.Font.Bold = $true

# 5) Copy the run-level language tag (en-US) from the immediately following
#    text ("Using cplex...") onto the new run -- Range.LanguageID isn't
#    wired to XML serialization in this host, so we clone formatting
#    (including the <w:lang> run property) via FormattedText instead, then
#    restore the intended text.
$langSample = $d.Range($pStart + 5, $pStart + 6)
$newRange.FormattedText = $langSample.FormattedText
$fixRange = $d.Range($pStart, $pStart + 1)
$fixRange.Text = " Tip "
